$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Swap the order of slides 4 ("Keep it simple stupid") and 5
#    ("Concurrent programming is necessary but dangerous in many languages")
# ---------------------------------------------------------------------------
$kiss = $p.Slides.Item(4)
$kiss.MoveTo(5)

# ---------------------------------------------------------------------------
# 2. Edit the "Keep it simple stupid" slide, now at position 5.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(5)
$body = $slide.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

# 2a. Insert a brand-new first bullet.
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.InsertBefore("A lot of the complexity and readability problems are from trying to package too much complexity together" + [char]13)

# 2b. Split "The first step to solving a complex problem is to break it into pieces"
#     into two runs: "The " + "first step to solving a complex problem is to break it into pieces"
$secondPara = $tr.Paragraphs(2, 1)
$theRun = $secondPara.Characters(1, 4)
$theRun.Text = $theRun.Text

# 2c. Split the last paragraph into three runs and add the trailing ellipsis.
$lastPara = $tr.Paragraphs(4, 1)
$midRun = $lastPara.Characters(84, 15)
$midRun.Text = $midRun.Text
$endRun = $lastPara.Characters(99, 10)
$endRun.Text = $endRun.Text + [char]0x2026

# 2d. Shrink text on overflow (adds <a:normAutofit/> to the placeholder).
$body.TextFrame.AutoSize = 2

# 2e. Extend the paragraph-click build to the newly-added first paragraph
#     (index 4, 1-based) so it animates in along with the rest.
$ms = $slide.TimeLine.MainSequence
$newEffect = $ms.AddEffect($body, 1)
$newEffect.Paragraph = 4
